# Adding bookmarks into the list 'Vek'
#
# 1) Reflow the 5 rounded-rectangle "age bucket" boxes on slide 2 to make
#    room for a 6th box, and add that 6th box.
# 2) Refresh the cached "today" date field on the slide master and every
#    slide layout (side effect of PowerPoint re-saving the deck on a
#    later day).

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# Left/Top/Width/Height are expressed in points; the literals below are the
# point value that round-trips to the exact target EMU (1 pt = 12700 EMU)
# through this host's float conversion.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Move the box that stays in place ("rohy 2") -------------------------
$rohy2 = Get-ShapeById $s.Shapes 3
$rohy2.Left   = 100.17393875187402
$rohy2.Top    = 325.35976377952755
$rohy2.Width  = 129.9575590551181
$rohy2.Height = 47.636614173228345

# --- Replace the remaining four boxes with narrower / repositioned ones --
# (duplicating "rohy 2" keeps the fill/line/shadow/style exactly the same)
$targets = @(
    @{ OldId = 8;  NewName = "Obdélník: se zakulacenými rohy 3"; X = 240.7792129524252;  Y = 326.0648956307638 },
    @{ OldId = 9;  NewName = "Obdélník: se zakulacenými rohy 4"; X = 381.3844881889764;  Y = 325.35976377952755 },
    @{ OldId = 10; NewName = "Obdélník: se zakulacenými rohy 5"; X = 521.9897766115276;  Y = 325.3596954353701 },
    @{ OldId = 12; NewName = "Obdélník: se zakulacenými rohy 6"; X = 662.5950622560788;  Y = 324.9799346928425 }
)

foreach ($t in $targets) {
    $old = Get-ShapeById $s.Shapes $t.OldId
    $old.Delete()
}

foreach ($t in $targets) {
    $new = $rohy2.Duplicate()
    $new.Name   = $t.NewName
    $new.Left   = $t.X
    $new.Top    = $t.Y
    $new.Width  = 129.9575590551181
    $new.Height = 47.636614173228345
}

# --- Add the brand-new sixth box ------------------------------------------
$new6 = $rohy2.Duplicate()
$new6.Name   = "Obdélník: se zakulacenými rohy 12"
$new6.Left   = 803.20031496063
$new6.Top    = 324.56661417322834
$new6.Width  = 129.9575590551181
$new6.Height = 47.636614173228345

# --- Refresh the cached date field (master + every layout) ---------------
$newDate = "12/27/2024"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
